# Add a "confidence" column (C) with predicted probability values,
# mirroring the style of the existing header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell C1 - match style of existing headers (A1/B1 use style index 1)
$ws.Range("C1").Value = "confidence"

$srcHeader = $ws.Range("A1")
$dstHeader = $ws.Range("C1")
$dstHeader.Font.Bold = $srcHeader.Font.Bold
$dstHeader.HorizontalAlignment = $srcHeader.HorizontalAlignment
$dstHeader.VerticalAlignment = $srcHeader.VerticalAlignment
$dstHeader.Borders.LineStyle = $srcHeader.Borders.LineStyle

# Confidence values for rows 2-21
$values = @{
    2  = 0.2867957072011674
    3  = 0.2984461658463237
    4  = 0.2261905999665737
    5  = 0.199830411225062
    6  = 0.2892984845204236
    7  = 0.199830411225062
    8  = 0.269856238839315
    9  = 0.1889446286839567
    10 = 0.2393715690698689
    11 = 0.2170463137454154
    12 = 0.2507262996613995
    13 = 0.199830411225062
    14 = 0.2228935537865555
    15 = 0.2892984845204236
    16 = 0.199830411225062
    17 = 0.2671192696684385
    18 = 0.273080436631102
    19 = 0.2039381202026201
    20 = 0.2301744125233183
    21 = 0.2456780320684599
}

foreach ($row in $values.Keys | Sort-Object) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
